$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 61; existing rows 61-190 shift down to 62-191.
$ws.Rows(61).Insert()

# Populate the newly inserted row 61 with the new data record.
$ws.Range("A61").Value = 10
$ws.Range("B61").Value = "Vega Modelo de Temuco"
$ws.Range("C61").Value = "La Araucanía"
$ws.Range("D61").Value = 44645
$ws.Range("E61").Value = 9
$ws.Range("F61").Value = 100112005
$ws.Range("G61").Value = "Puerro"
$ws.Range("H61").Value = "Azul de Maquehue"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 45
$ws.Range("K61").Value = 12000
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = 12000
$ws.Range("N61").Value = "$/docena de paquetes"
$ws.Range("O61").Value = "Provincia de Cautín"
$ws.Range("P61").Value = 1000
$ws.Range("Q61").Value = 12
$ws.Range("R61").Value = "Hortaliza"
